# The document contains three `<div>` sections, each starting with an
# `<id>...</id>` tag that was originally split across three separate
# runs (one for the literal "<id>", one for the Arial-styled id value,
# and one for the literal "</id>"). This collapses each of those into
# a single run (using the formatting of the first/opening-tag run:
# Courier New, color 7f6000, sz/szCs 18) that carries the whole
# "<id>...</id>" string as one piece of text.

$d = $word.ActiveDocument

$ids = @("p068r_1", "p068r_2", "p068r_3")

foreach ($id in $ids) {
    $needle = "<id>$id</id>"
    $r = $d.Content
    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
    # MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace(2=All)
    $r.Find.Execute($needle, $true, $false, $false, $false, $false, `
                     $true, 1, $false, $needle, 2) | Out-Null
}
